# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-236) from 2023-09-23 (serial 45192) to 2023-10-03 (serial 45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C236").Value = 45202
